$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "A120"
$ws.Range("D6").Value = "B20"
$ws.Range("D7").Value = "C20"

$ws.Range("A8").Value = "DV DSS Path Indices"
$ws.Range("B8").Value = "DV Dss path names"
$ws.Range("D8").Value = "G20"

$ws.Range("A9").Value = "SV DSS Path Indices"
$ws.Range("B9").Value = "SV Dss path names"
$ws.Range("D9").Value = "H20"

$ws.Range("D10").Value = "I20"
$ws.Range("D11").Value = "J20"

$ws.Range("A5:D11").Select()
